# Adds 20 new trading-day rows (631-650) to the Luxury_Goods NAV sheet,
# re-appending the 2024-09-02 .. 2024-09-27 date range with refreshed
# price / basket-value / return / NAV figures, and extends the sheet's
# used-range dimension from J630 to J650 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column layout (row 1 header): A=Date, C..G=constituent prices,
# H=Basket Value, I=daily return, J=NAV. Column B ("Quantities") is only
# populated on row 2, so new data rows leave it blank.
$cols = @("A","C","D","E","F","G","H","I","J")

# Each inner array: Date, C, D, E, F, G, H, I, J
$newRows = @(
    @("2024-09-02",1660.599975585938,1280.449951171875,1769.650024414062,2220.5,464.8500061035156,8325.749969482422,0,229.6222988627398),
    @("2024-09-03",1683.75,1248.449951171875,1810.949951171875,2205.5,469.6000061035156,8357.449920654297,0.003807458942205738,230.4965763378746),
    @("2024-09-04",1699.5,1288.900024414062,1780.25,2217.949951171875,466.5499877929688,8386.249938964844,0.003446029420932761,231.2908743213592),
    @("2024-09-05",1694.699951171875,1272.300048828125,1783.150024414062,2236.14990234375,467.75,8389.549926757812,0.0003934998142180441,231.381887237435),
    @("2024-09-06",1724.449951171875,1234.300048828125,1778.650024414062,2224.39990234375,463.3999938964844,8351.999908447266,-0.004475808432915339,230.3462662353138),
    @("2024-09-09",1789.300048828125,1234.300048828125,1817.949951171875,2138.60009765625,460.7999877929688,8362.550109863281,0.001263194627833399,230.6372384013637),
    @("2024-09-10",1799.949951171875,1249.650024414062,1830.099975585938,2140.14990234375,466.9500122070312,8420.699890136719,0.006953594239734628,232.2409961737797),
    @("2024-09-11",1788.25,1255,1828.699951171875,2119.85009765625,459.2999877929688,8369.700012207031,-0.006056489198650145,230.8344310889695),
    @("2024-09-12",1816.650024414062,1257.449951171875,1852.949951171875,2111,459,8415.049926757812,0.005418344084571652,232.0851714631759),
    @("2024-09-13",1888,1264.349975585938,1921.550048828125,2101.35009765625,451.75,8530.500122070312,0.01371949023681921,235.2692617071754),
    @("2024-09-16",1934.900024414062,1259.75,1916,2047.199951171875,447.8999938964844,8501.549957275391,-0.003393724210849177,234.4708227176512),
    @("2024-09-17",1891.199951171875,1257.550048828125,1904.050048828125,2034.699951171875,452.75,8445.75,-0.006563504014657771,232.9318725314237),
    @("2024-09-18",1879.449951171875,1232.050048828125,1889.400024414062,2042.550048828125,448.8500061035156,8390.000091552734,-0.006600942302017657,231.3943026805429),
    @("2024-09-19",1866.650024414062,1234.5,1929.199951171875,2025.699951171875,450,8406.049926757812,0.00191297199403341,231.8369535011497),
    @("2024-09-20",1930.099975585938,1246.550048828125,1928.400024414062,2062.300048828125,446.7999877929688,8507.750061035156,0.01209844518691423,234.6418201753845),
    @("2024-09-23",2012.849975585938,1229.900024414062,1882.449951171875,2039,445.75,8501.449951171875,-0.000740514215636784,234.4680645719617),
    @("2024-09-24",2026,1238.800048828125,1914.400024414062,2064.550048828125,439.25,8561.500122070312,0.007063521075032611,236.1242346874879),
    @("2024-09-25",2051.39990234375,1252.949951171875,1910.849975585938,2112.050048828125,441.5499877929688,8651.899841308594,0.01055886444540762,238.6174384738287),
    @("2024-09-26",2134.800048828125,1277.199951171875,1857.550048828125,2155.60009765625,419.25,8682.900146484375,0.003583063343818423,239.4724198708201),
    @("2024-09-27",2119.5,1275.75,1866.699951171875,2201.449951171875,422.7999877929688,8731.799865722656,0.005631726544509473,240.8210630544845)
)

$startRow = 631
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $newRows[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $col = $cols[$j]
        $cellRef = "$col$r"
        if ($j -eq 0) {
            # Force the date to be written as literal text (matching the
            # existing A-column cells), not auto-converted to a date serial.
            $ws.Range($cellRef).NumberFormat = "@"
        }
        $ws.Range($cellRef).Value = $rowVals[$j]
    }
}
